# This workbook tracks daily "Puerro" (leek) price records for the
# "Vega Modelo de Temuco" market. The edit inserts one new daily record
# (a new row) above the current row 29, which pushes all the existing
# records down by one row (the former last row, 122, becomes row 123).
#
# The newly inserted row reuses the same market/category/variety/quality
# metadata as the record that used to sit at row 29, but carries its own
# date (column D) and volume (column J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 29; this shifts rows 29:122 down to
# 30:123 and keeps all their existing data/styles intact.
$ws.Rows(29).Insert()

# Populate the newly inserted row 29 with the new record.
$ws.Range("A29").Value = 10
$ws.Range("B29").Value = "Vega Modelo de Temuco"
$ws.Range("C29").Value = "La Araucanía"
$ws.Range("D29").Value = 44453
$ws.Range("E29").Value = 9
$ws.Range("F29").Value = 100112005
$ws.Range("G29").Value = "Puerro"
$ws.Range("H29").Value = "Azul de Maquehue"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 40
$ws.Range("K29").Value = 8000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 8000
$ws.Range("N29").Value = "$/docena de paquetes"
$ws.Range("O29").Value = "Provincia de Cautín"
$ws.Range("P29").Value = 667
$ws.Range("Q29").Value = 12
$ws.Range("R29").Value = "Hortaliza"
